$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header row: "_old" -> "_FV2410", "_new" -> "_FV2504" ---
# Columns A..J (1..10) currently end with "_old" (except "diff" in K/11 which is unchanged)
for ($c = 1; $c -le 10; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $val = $cell.Value()
    if ($val -like "*_old") {
        $cell.Value = ($val -replace "_old$", "_FV2410")
    }
}

# Columns L..U (12..21) currently end with "_new"
for ($c = 12; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $val = $cell.Value()
    if ($val -like "*_new") {
        $cell.Value = ($val -replace "_new$", "_FV2504")
    }
}

# --- 2. Turn the data range into an Excel Table (ListObject) ---
$range = $ws.Range("A1:U80")
$lo = $ws.ListObjects.Add(1, $range, $null, 1)
$lo.Name = "Table1"
# No named/banded table style in the target workbook - clear it
$lo.TableStyle = $null
$lo.ShowTableStyleRowStripes = $true
$lo.ShowTableStyleColumnStripes = $false

# --- 3. Freeze the header row (split after row 1) ---
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$win = $excel.ActiveWindow
$win.FreezePanes = $true
